$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-obsolete extra job rows (6-9), collapsing the report to
# just two rows of data (rows 4 and 5).
$ws.Rows("6:9").Delete()

# Re-anchor the conditional formatting to the now-smaller data range.
$fcs = $ws.Range("A4:G6").FormatConditions
for ($i = 1; $i -le $fcs.Count; $i++) {
    $fcs.Item($i).ModifyAppliesToRange($ws.Range("A4:G6"))
}

# Report title / timestamp.
$ws.Range("B1").Value = "Relatório das Ultimas 24 horas 2022-09-26 14:09:21"

# Header rename: "Task Code" -> "Job Code".
$ws.Range("B3").Value = "Job Code"

# Row 4: first job entry.
$ws.Range("B4").Value = "01-job"
$ws.Range("C4").Value = "CurrentCurrencyTrades"
$ws.Range("D4").Value = "Current exchange rates update"
$ws.Range("E4").Value = "Yes"
$ws.Range("G4").Value = "2022-09-26 13:45:59"

# Row 5: second job entry (now the same job, later run).
$ws.Range("B5").Value = "01-job"
$ws.Range("C5").Value = "CurrentCurrencyTrades"
$ws.Range("D5").Value = "Current exchange rates update"
$ws.Range("E5").Value = "Yes"
$ws.Range("G5").Value = "2022-09-26 13:46:28"
